# Simulink "Divsalar and Simon Compare" workbook update
# Adds an "N=3" data series (modulo-step improvement) alongside the
# existing "N=2" (renamed from "My Model" to "P. N=2") data, for all
# three SNR tables (DBPSK, DQPSK, 8-PSK) and updates the three charts
# to reflect the new series / titles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Worksheet data -- header labels (order matters: write the "P. N=3"
#    labels first so the shared-string table orders them the same way
#    the original authoring session did).
# ---------------------------------------------------------------------
$ws.Range("D2").Value = "P. N=3"
$ws.Range("C2").Value = "P. N=2"

$ws.Range("M2").Value = "P. N=3"
$ws.Range("L2").Value = "P. N=2"

$ws.Range("S2").Value = "P. N=3"
$ws.Range("R2").Value = "P. N=2"

# ---------------------------------------------------------------------
# 2. Worksheet data -- new "N=3" BER columns
# ---------------------------------------------------------------------

# DBPSK (M=2), SNR 7..10 in A3:A6
$ws.Range("D3").Formula = "=0.00307"
$ws.Range("D4").Formula = "=88/100000"
$ws.Range("D5").Value = 0.00016
$ws.Range("D5").NumberFormat = "0.00E+00"
$ws.Range("D6").Formula = "=0.000017"
$ws.Range("D6").NumberFormat = "0.00E+00"

# DQPSK (M=4), SNR 7..12 in J3:J8
$ws.Range("M3").Formula = "=0.01379"
$ws.Range("M4").Formula = "=0.00576"
$ws.Range("M5").Formula = "=0.00196"
$ws.Range("M6").Formula = "=0.00049"
$ws.Range("M7").Formula = "=0.000076"
$ws.Range("M7").NumberFormat = "0.00E+00"
$ws.Range("M8").Formula = "=0.000012"

# 8-PSK (M=8), SNR 11..16 in P3:P8
$ws.Range("S3").Formula = "=0.01232"
$ws.Range("S4").Formula = "=0.004618"
$ws.Range("S5").Formula = "=0.001342"
$ws.Range("S6").Formula = "=0.000292"
$ws.Range("S7").Formula = "=0.000034"
$ws.Range("S8").Formula = "=0.000008"

# ---------------------------------------------------------------------
# 3. Selection state, matching the saved session.
# ---------------------------------------------------------------------
$ws.Range("U6").Select()

# ---------------------------------------------------------------------
# 4. Charts: update titles ("N=2" -> "N=2,3") and series.
# ---------------------------------------------------------------------
$chartObjs = $ws.ChartObjects()

# --- Chart 6 (rId1 / chart1.xml): M=2, N=2 scatter -------------------
$co1 = $chartObjs.Item(1)
$chart1 = $co1.Chart
$chart1.ChartTitle.Text = "Non-Coherent Detection M=2, N=2,3"

$sc1 = $chart1.SeriesCollection()
# Drop the two duplicate ("marker-less") series, keeping one Div&Sim
# series and one My-Model series to rename/re-use.
$sc1.Item(3).Delete()
$sc1b = $chart1.SeriesCollection()
$sc1b.Item(2).Delete()
$sc1c = $chart1.SeriesCollection()
$sc1c.Item(1).Name = "Div&Simon"
$sc1c.Item(2).Name = "P. N=2"

$newSer1 = $chart1.SeriesCollection().NewSeries()
$newSer1.Name = "P. N=3"
$newSer1.XValues = $ws.Range("A3:A6")
$newSer1.Values = $ws.Range("D3:D6")

# --- Chart 1 (rId2 / chart2.xml): M=4, N=2 scatter --------------------
$co2 = $chartObjs.Item(2)
$chart2 = $co2.Chart
$chart2.ChartTitle.Text = "Non-Coherent Detection M=4, N=2,3"

$sc2 = $chart2.SeriesCollection()
$sc2.Item(2).Name = "P. N=2"

$newSer2 = $chart2.SeriesCollection().NewSeries()
$newSer2.Name = "N=3 Peyton"
$newSer2.XValues = $ws.Range("J3:J8")
$newSer2.Values = $ws.Range("M3:M8")

# --- Chart 7 (rId3 / chart3.xml): M=8, N=2 scatter --------------------
$co3 = $chartObjs.Item(3)
$chart3 = $co3.Chart
$chart3.ChartTitle.Text = "Non-Coherent Detection M=8, N=2,3"

$sc3 = $chart3.SeriesCollection()
$sc3.Item(2).Name = "P. N=2"

$newSer3 = $chart3.SeriesCollection().NewSeries()
$newSer3.Name = "P. N=3"
$newSer3.XValues = $ws.Range("P3:P8")
$newSer3.Values = $ws.Range("S3:S8")
